# Add Datas of prison-copes & doctors in XML
# Target sheet: Sheet1 ("Charactor" XML-mapped table), rows 1:E113 -> 1:E118
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Grow the XML-mapped table / sheet range from A1:E113 to A1:E118
#    so the new rows inherit the Text ("@") number format used by the
#    rest of the mapped data (style index 1 in the original file).
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$ws.Range("A114:E118").NumberFormat = "@"
$lo.Resize($ws.Range("A1:E118"))

# ------------------------------------------------------------------
# 2. Rewrite rows 2-6 : BodyType/SkinColor columns move from the old
#    fat/thin/strong + normal/brown/grey/white scheme to the new
#    1/2/3 scheme, and SkinColor (D) is cleared.
# ------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = "GEEK"
$ws.Cells.Item(2,2).Value = "1"
$ws.Cells.Item(2,3).Value = "1"
$ws.Cells.Item(2,4).Value = ""
$ws.Cells.Item(2,5).Value = "Body_up"

$ws.Cells.Item(3,1).Value = "SMARTY"
$ws.Cells.Item(3,2).Value = "2"
$ws.Cells.Item(3,3).Value = "2"
$ws.Cells.Item(3,4).Value = ""
$ws.Cells.Item(3,5).Value = "Leg_up"

$ws.Cells.Item(4,1).Value = "RED"
$ws.Cells.Item(4,2).Value = "3"
$ws.Cells.Item(4,3).Value = "2"
$ws.Cells.Item(4,4).Value = ""
$ws.Cells.Item(4,5).Value = "Leg_up"

$ws.Cells.Item(5,1).Value = "RANBOW"
$ws.Cells.Item(5,2).Value = "4"
$ws.Cells.Item(5,3).Value = "3"
$ws.Cells.Item(5,4).Value = ""
$ws.Cells.Item(5,5).Value = "Leg_up"

$ws.Cells.Item(6,1).Value = "OLD FISH"
$ws.Cells.Item(6,2).Value = "5"
$ws.Cells.Item(6,3).Value = "2"
$ws.Cells.Item(6,4).Value = ""
$ws.Cells.Item(6,5).Value = "Leg_up"

# ------------------------------------------------------------------
# 3. Row 104 only changes its HeadNO (B) value.
# ------------------------------------------------------------------
$ws.Cells.Item(104,2).Value = "102_2"

# ------------------------------------------------------------------
# 4. Rows 105-113 become the new prison-guard / doctor records
#    (previously only had a HeadNO in column B).
# ------------------------------------------------------------------
$ws.Cells.Item(105,1).Value = "DOCTOR"
$ws.Cells.Item(105,2).Value = "doc"
$ws.Cells.Item(105,3).Value = "doctor"
$ws.Cells.Item(105,4).Value = ""
$ws.Cells.Item(105,5).Value = "Body_up"

$ws.Cells.Item(106,1).Value = "BELLY"
$ws.Cells.Item(106,2).Value = "g1"
$ws.Cells.Item(106,3).Value = "g"
$ws.Cells.Item(106,4).Value = ""
$ws.Cells.Item(106,5).Value = "Body_up"

$ws.Cells.Item(107,1).Value = "MR SLAVE"
$ws.Cells.Item(107,2).Value = "g2"
$ws.Cells.Item(107,3).Value = "g"
$ws.Cells.Item(107,4).Value = ""
$ws.Cells.Item(107,5).Value = "Body_up"

$ws.Cells.Item(108,1).Value = "VIC MICKEY"
$ws.Cells.Item(108,2).Value = "g3"
$ws.Cells.Item(108,3).Value = "g"
$ws.Cells.Item(108,4).Value = ""
$ws.Cells.Item(108,5).Value = "Body_up"

$ws.Cells.Item(109,1).Value = "SHANE"
$ws.Cells.Item(109,2).Value = "g4"
$ws.Cells.Item(109,3).Value = "g"
$ws.Cells.Item(109,4).Value = ""
$ws.Cells.Item(109,5).Value = "Body_up"

$ws.Cells.Item(110,1).Value = "GORDON"
$ws.Cells.Item(110,2).Value = "g5"
$ws.Cells.Item(110,3).Value = "g"
$ws.Cells.Item(110,4).Value = ""
$ws.Cells.Item(110,5).Value = "Body_up"

$ws.Cells.Item(111,1).Value = "PIKMAL"
$ws.Cells.Item(111,2).Value = "g6"
$ws.Cells.Item(111,3).Value = "g"
$ws.Cells.Item(111,4).Value = ""
$ws.Cells.Item(111,5).Value = "Body_up"

$ws.Cells.Item(112,1).Value = "DOCTOR SNLAD"
$ws.Cells.Item(112,2).Value = "g10"
$ws.Cells.Item(112,3).Value = "doctor"
$ws.Cells.Item(112,4).Value = ""
$ws.Cells.Item(112,5).Value = "Body_up"

$ws.Cells.Item(113,1).Value = "DOCTOR SNLAD"
$ws.Cells.Item(113,2).Value = "g10_"
$ws.Cells.Item(113,3).Value = "doctor"
$ws.Cells.Item(113,4).Value = ""
$ws.Cells.Item(113,5).Value = "Body_up"

# ------------------------------------------------------------------
# 5. Brand new rows 114-118.
# ------------------------------------------------------------------
$ws.Cells.Item(114,2).Value = "g11"

$ws.Cells.Item(115,2).Value = "g11_kidnap"

$ws.Cells.Item(116,1).Value = "MAD DOCTOR"
$ws.Cells.Item(116,2).Value = "g12"
$ws.Cells.Item(116,3).Value = "doctor"
$ws.Cells.Item(116,4).Value = ""
$ws.Cells.Item(116,5).Value = "Body_up"

$ws.Cells.Item(117,2).Value = "g13"

$ws.Cells.Item(118,1).Value = "COOK"
$ws.Cells.Item(118,2).Value = "g14"
$ws.Cells.Item(118,3).Value = "1"
$ws.Cells.Item(118,4).Value = ""
$ws.Cells.Item(118,5).Value = "Body_up"

# ------------------------------------------------------------------
# 6. Column A widens a bit to fit the new longer names.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.66

# ------------------------------------------------------------------
# 7. Data validation lists are refreshed for the new BodyType /
#    SkinColor vocab, and a couple of new ranges receive the same
#    SkinColor rule reserved for rows above/below the table body.
# ------------------------------------------------------------------
$ws.Range("D1:D1048576").Validation.Delete()
$ws.Range("D1").Validation.Add(3, 1, 1, """_b,_g,_w,,""")
$ws.Range("D119:D1048576").Validation.Add(3, 1, 1, """_b,_g,_w,,""")

$ws.Range("C1:C1048576").Validation.Delete()
$ws.Range("C1:C1048576").Validation.Add(3, 1, 1, """1,2,3,g,doctor""")

$ws.Range("D2:D118").Validation.Add(3, 1, 1, """_b,_g,_w""")

# ------------------------------------------------------------------
# 8. Selection / scroll position, matching the author's view when
#    they finished entering the new rows.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 106
$ws.Range("B2:B118").Select()
